# Title-page revisions (R2):
#  1. Colour the word "Only" (in the title) blue (RGB 0,176,240 / hex 00B0F0),
#     leaving the following space in its own, uncoloured run.
#  2. Highlight the "XXX" word-count placeholder in yellow.

$d = $word.ActiveDocument

# --- 1. Colour "Only" -------------------------------------------------
$rngOnly = $d.Content
$foundOnly = $rngOnly.Find.Execute("Only", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if ($foundOnly) {
    # Setting the font colour on just this sub-range causes Word to split
    # the run, so the trailing space keeps its original (uncoloured) run.
    $rngOnly.Font.Color = 15773696   # RGB(0,176,240) -> wdColor 0x00F0B000
}

# --- 2. Highlight "XXX" -----------------------------------------------
# NOTE: setting Range.HighlightColorIndex directly on a sub-range mutates
# every run in the enclosing paragraph in this host, not just the target
# run. Driving the highlight through Find/Replace's replacement formatting
# keeps it scoped to only the matched text ("XXX"), which is what we want.
$rngXXX = $d.Content
$rngXXX.Find.ClearFormatting()
$rngXXX.Find.Replacement.ClearFormatting()
$rngXXX.Find.Replacement.Highlight = $true
$rngXXX.Find.Execute("XXX", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "XXX", 1) | Out-Null
